# Update "想去人数" (number of people interested) counts in the
# "展览" and "全部类型" worksheets.
#
# 展览 (sheet1) uses data rows 2-10 -> column F
# 全部类型 (sheet4) uses data rows 2-11 -> column F (row 4 is a
#   "演出" entry that is unchanged, so it is skipped below)

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAll        = $wb.Worksheets.Item("全部类型")

# Sheet "展览": row -> new value
$exhibitionUpdates = @{
    2  = 741
    3  = 41
    4  = 251
    5  = 3039
    6  = 61
    7  = 3848
    9  = 966
    10 = 33
}

foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型": row -> new value
$allTypesUpdates = @{
    2  = 741
    3  = 41
    5  = 251
    6  = 3039
    7  = 61
    8  = 3848
    10 = 966
    11 = 33
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allTypesUpdates[$row]
}

$wb.Save()
